$d = $word.ActiveDocument

$d.Content.Find.Execute("99-54=45", $true, $false, $false, $false, $false, $true, 1, $false, "77+21=98", 2) | Out-Null
$d.Content.Find.Execute("70+4=74", $true, $false, $false, $false, $false, $true, 1, $false, "26-2=24", 2) | Out-Null
$d.Content.Find.Execute("55-39=16", $true, $false, $false, $false, $false, $true, 1, $false, "0+82=82", 2) | Out-Null
$d.Content.Find.Execute("28+68=96", $true, $false, $false, $false, $false, $true, 1, $false, "52+28=80", 2) | Out-Null
$d.Content.Find.Execute("38+12=50", $true, $false, $false, $false, $false, $true, 1, $false, "45+54=99", 2) | Out-Null
$d.Content.Find.Execute("70-38=32", $true, $false, $false, $false, $false, $true, 1, $false, "1+41=42", 2) | Out-Null
$d.Content.Find.Execute("50-19=31", $true, $false, $false, $false, $false, $true, 1, $false, "86-1=85", 2) | Out-Null
$d.Content.Find.Execute("40+41=81", $true, $false, $false, $false, $false, $true, 1, $false, "75-2=73", 2) | Out-Null
$d.Content.Find.Execute("93-14=79", $true, $false, $false, $false, $false, $true, 1, $false, "22+6=28", 2) | Out-Null
$d.Content.Find.Execute("80-27=53", $true, $false, $false, $false, $false, $true, 1, $false, "41+16=57", 2) | Out-Null
$d.Content.Find.Execute("97-64=33", $true, $false, $false, $false, $false, $true, 1, $false, "44+26=70", 2) | Out-Null
$d.Content.Find.Execute("80-76=4", $true, $false, $false, $false, $false, $true, 1, $false, "30-10=20", 2) | Out-Null
$d.Content.Find.Execute("48-19=29", $true, $false, $false, $false, $false, $true, 1, $false, "78+13=91", 2) | Out-Null
$d.Content.Find.Execute("70-60=10", $true, $false, $false, $false, $false, $true, 1, $false, "93-40=53", 2) | Out-Null
$d.Content.Find.Execute("6+74=80", $true, $false, $false, $false, $false, $true, 1, $false, "87-82=5", 2) | Out-Null
$d.Content.Find.Execute("73+17=90", $true, $false, $false, $false, $false, $true, 1, $false, "55-31=24", 2) | Out-Null
$d.Content.Find.Execute("0+66=66", $true, $false, $false, $false, $false, $true, 1, $false, "49+39=88", 2) | Out-Null
$d.Content.Find.Execute("92-91=1", $true, $false, $false, $false, $false, $true, 1, $false, "35+64=99", 2) | Out-Null
$d.Content.Find.Execute("13+85=98", $true, $false, $false, $false, $false, $true, 1, $false, "44+36=80", 2) | Out-Null
$d.Content.Find.Execute("28+16=44", $true, $false, $false, $false, $false, $true, 1, $false, "70-17=53", 2) | Out-Null
$d.Content.Find.Execute("14+3=17", $true, $false, $false, $false, $false, $true, 1, $false, "10-2=8", 2) | Out-Null
$d.Content.Find.Execute("43+1=44", $true, $false, $false, $false, $false, $true, 1, $false, "44+30=74", 2) | Out-Null
$d.Content.Find.Execute("66+14=80", $true, $false, $false, $false, $false, $true, 1, $false, "73-43=30", 2) | Out-Null
$d.Content.Find.Execute("99-88=11", $true, $false, $false, $false, $false, $true, 1, $false, "68-2=66", 2) | Out-Null
$d.Content.Find.Execute("37+2=39", $true, $false, $false, $false, $false, $true, 1, $false, "42+44=86", 2) | Out-Null
$d.Content.Find.Execute("71-10=61", $true, $false, $false, $false, $false, $true, 1, $false, "77+2=79", 2) | Out-Null
$d.Content.Find.Execute("15+33=48", $true, $false, $false, $false, $false, $true, 1, $false, "25+6=31", 2) | Out-Null
$d.Content.Find.Execute("6+81=87", $true, $false, $false, $false, $false, $true, 1, $false, "65-44=21", 2) | Out-Null
$d.Content.Find.Execute("16+63=79", $true, $false, $false, $false, $false, $true, 1, $false, "36+38=74", 2) | Out-Null
$d.Content.Find.Execute("84-48=36", $true, $false, $false, $false, $false, $true, 1, $false, "80+9=89", 2) | Out-Null
$d.Content.Find.Execute("46+6=52", $true, $false, $false, $false, $false, $true, 1, $false, "89-80=9", 2) | Out-Null
$d.Content.Find.Execute("71+15=86", $true, $false, $false, $false, $false, $true, 1, $false, "3+53=56", 2) | Out-Null
$d.Content.Find.Execute("74+10=84", $true, $false, $false, $false, $false, $true, 1, $false, "27+42=69", 2) | Out-Null
$d.Content.Find.Execute("31-10=21", $true, $false, $false, $false, $false, $true, 1, $false, "54+33=87", 2) | Out-Null
$d.Content.Find.Execute("93-60=33", $true, $false, $false, $false, $false, $true, 1, $false, "61-20=41", 2) | Out-Null
$d.Content.Find.Execute("24-5=19", $true, $false, $false, $false, $false, $true, 1, $false, "61-59=2", 2) | Out-Null
$d.Content.Find.Execute("58+2=60", $true, $false, $false, $false, $false, $true, 1, $false, "45-22=23", 2) | Out-Null
$d.Content.Find.Execute("65-57=8", $true, $false, $false, $false, $false, $true, 1, $false, "1+73=74", 2) | Out-Null
$d.Content.Find.Execute("62-14=48", $true, $false, $false, $false, $false, $true, 1, $false, "52-28=24", 2) | Out-Null
$d.Content.Find.Execute("96-82=14", $true, $false, $false, $false, $false, $true, 1, $false, "7+2=9", 2) | Out-Null
$d.Content.Find.Execute("17+74=91", $true, $false, $false, $false, $false, $true, 1, $false, "41+56=97", 2) | Out-Null
$d.Content.Find.Execute("23+51=74", $true, $false, $false, $false, $false, $true, 1, $false, "12+15=27", 2) | Out-Null
$d.Content.Find.Execute("34+41=75", $true, $false, $false, $false, $false, $true, 1, $false, "99-42=57", 2) | Out-Null
$d.Content.Find.Execute("89+7=96", $true, $false, $false, $false, $false, $true, 1, $false, "53+34=87", 2) | Out-Null
$d.Content.Find.Execute("83-55=28", $true, $false, $false, $false, $false, $true, 1, $false, "38+32=70", 2) | Out-Null
$d.Content.Find.Execute("60-58=2", $true, $false, $false, $false, $false, $true, 1, $false, "28+17=45", 2) | Out-Null
$d.Content.Find.Execute("14+16=30", $true, $false, $false, $false, $false, $true, 1, $false, "57-15=42", 2) | Out-Null
$d.Content.Find.Execute("43-6=37", $true, $false, $false, $false, $false, $true, 1, $false, "91-38=53", 2) | Out-Null
$d.Content.Find.Execute("16+40=56", $true, $false, $false, $false, $false, $true, 1, $false, "49-27=22", 2) | Out-Null
$d.Content.Find.Execute("62-38=24", $true, $false, $false, $false, $false, $true, 1, $false, "99-29=70", 2) | Out-Null
$d.Content.Find.Execute("99-8=91", $true, $false, $false, $false, $false, $true, 1, $false, "53-13=40", 2) | Out-Null
$d.Content.Find.Execute("37+13=50", $true, $false, $false, $false, $false, $true, 1, $false, "45-20=25", 2) | Out-Null
$d.Content.Find.Execute("28+9=37", $true, $false, $false, $false, $false, $true, 1, $false, "55+6=61", 2) | Out-Null
$d.Content.Find.Execute("82-41=41", $true, $false, $false, $false, $false, $true, 1, $false, "1+39=40", 2) | Out-Null
$d.Content.Find.Execute("67-18=49", $true, $false, $false, $false, $false, $true, 1, $false, "73+5=78", 2) | Out-Null
$d.Content.Find.Execute("64+32=96", $true, $false, $false, $false, $false, $true, 1, $false, "88+4=92", 2) | Out-Null
$d.Content.Find.Execute("35+1=36", $true, $false, $false, $false, $false, $true, 1, $false, "95-7=88", 2) | Out-Null
$d.Content.Find.Execute("29-8=21", $true, $false, $false, $false, $false, $true, 1, $false, "14+42=56", 2) | Out-Null
$d.Content.Find.Execute("49+40=89", $true, $false, $false, $false, $false, $true, 1, $false, "14+6=20", 2) | Out-Null
$d.Content.Find.Execute("33+0=33", $true, $false, $false, $false, $false, $true, 1, $false, "32-23=9", 2) | Out-Null
$d.Content.Find.Execute("8+8=16", $true, $false, $false, $false, $false, $true, 1, $false, "60-6=54", 2) | Out-Null
$d.Content.Find.Execute("28-15=13", $true, $false, $false, $false, $false, $true, 1, $false, "63+6=69", 2) | Out-Null
$d.Content.Find.Execute("1+10=11", $true, $false, $false, $false, $false, $true, 1, $false, "60-7=53", 2) | Out-Null
$d.Content.Find.Execute("17+65=82", $true, $false, $false, $false, $false, $true, 1, $false, "42+22=64", 2) | Out-Null
$d.Content.Find.Execute("76-52=24", $true, $false, $false, $false, $false, $true, 1, $false, "43-18=25", 2) | Out-Null
$d.Content.Find.Execute("21-18=3", $true, $false, $false, $false, $false, $true, 1, $false, "34+3=37", 2) | Out-Null
$d.Content.Find.Execute("30-24=6", $true, $false, $false, $false, $false, $true, 1, $false, "56+10=66", 2) | Out-Null
$d.Content.Find.Execute("53-17=36", $true, $false, $false, $false, $false, $true, 1, $false, "42+19=61", 2) | Out-Null
$d.Content.Find.Execute("30+2=32", $true, $false, $false, $false, $false, $true, 1, $false, "56-24=32", 2) | Out-Null
$d.Content.Find.Execute("15+55=70", $true, $false, $false, $false, $false, $true, 1, $false, "20-17=3", 2) | Out-Null
$d.Content.Find.Execute("75-45=30", $true, $false, $false, $false, $false, $true, 1, $false, "77+12=89", 2) | Out-Null
$d.Content.Find.Execute("6+62=68", $true, $false, $false, $false, $false, $true, 1, $false, "1+59=60", 2) | Out-Null
$d.Content.Find.Execute("27+26=53", $true, $false, $false, $false, $false, $true, 1, $false, "19+47=66", 2) | Out-Null
$d.Content.Find.Execute("47+46=93", $true, $false, $false, $false, $false, $true, 1, $false, "83-73=10", 2) | Out-Null
$d.Content.Find.Execute("9+2=11", $true, $false, $false, $false, $false, $true, 1, $false, "35-5=30", 2) | Out-Null
$d.Content.Find.Execute("22+70=92", $true, $false, $false, $false, $false, $true, 1, $false, "19+5=24", 2) | Out-Null
$d.Content.Find.Execute("18+63=81", $true, $false, $false, $false, $false, $true, 1, $false, "62-17=45", 2) | Out-Null
$d.Content.Find.Execute("31-22=9", $true, $false, $false, $false, $false, $true, 1, $false, "40-8=32", 2) | Out-Null
$d.Content.Find.Execute("32+23=55", $true, $false, $false, $false, $false, $true, 1, $false, "33+57=90", 2) | Out-Null
$d.Content.Find.Execute("25+20=45", $true, $false, $false, $false, $false, $true, 1, $false, "21-5=16", 2) | Out-Null
$d.Content.Find.Execute("2+42=44", $true, $false, $false, $false, $false, $true, 1, $false, "61-52=9", 2) | Out-Null
$d.Content.Find.Execute("49-20=29", $true, $false, $false, $false, $false, $true, 1, $false, "91-39=52", 2) | Out-Null
$d.Content.Find.Execute("82-35=47", $true, $false, $false, $false, $false, $true, 1, $false, "67+8=75", 2) | Out-Null
$d.Content.Find.Execute("78-27=51", $true, $false, $false, $false, $false, $true, 1, $false, "44+45=89", 2) | Out-Null
$d.Content.Find.Execute("84-30=54", $true, $false, $false, $false, $false, $true, 1, $false, "18+14=32", 2) | Out-Null
$d.Content.Find.Execute("3+96=99", $true, $false, $false, $false, $false, $true, 1, $false, "21+46=67", 2) | Out-Null
$d.Content.Find.Execute("87-69=18", $true, $false, $false, $false, $false, $true, 1, $false, "20+25=45", 2) | Out-Null
$d.Content.Find.Execute("58-24=34", $true, $false, $false, $false, $false, $true, 1, $false, "37-21=16", 2) | Out-Null
$d.Content.Find.Execute("39+33=72", $true, $false, $false, $false, $false, $true, 1, $false, "45-44=1", 2) | Out-Null
$d.Content.Find.Execute("76-42=34", $true, $false, $false, $false, $false, $true, 1, $false, "15-3=12", 2) | Out-Null
$d.Content.Find.Execute("58+36=94", $true, $false, $false, $false, $false, $true, 1, $false, "28-27=1", 2) | Out-Null
$d.Content.Find.Execute("22+71=93", $true, $false, $false, $false, $false, $true, 1, $false, "11+34=45", 2) | Out-Null
$d.Content.Find.Execute("75+19=94", $true, $false, $false, $false, $false, $true, 1, $false, "69-42=27", 2) | Out-Null
$d.Content.Find.Execute("27-10=17", $true, $false, $false, $false, $false, $true, 1, $false, "15+59=74", 2) | Out-Null
$d.Content.Find.Execute("81-30=51", $true, $false, $false, $false, $false, $true, 1, $false, "50-3=47", 2) | Out-Null
$d.Content.Find.Execute("89-77=12", $true, $false, $false, $false, $false, $true, 1, $false, "37+23=60", 2) | Out-Null
$d.Content.Find.Execute("63+25=88", $true, $false, $false, $false, $false, $true, 1, $false, "13+8=21", 2) | Out-Null
$d.Content.Find.Execute("51-29=22", $true, $false, $false, $false, $false, $true, 1, $false, "44+31=75", 2) | Out-Null
$d.Content.Find.Execute("97-32=65", $true, $false, $false, $false, $false, $true, 1, $false, "33+12=45", 2) | Out-Null
$d.Content.Find.Execute("16+77=93", $true, $false, $false, $false, $false, $true, 1, $false, "58-40=18", 2) | Out-Null
